$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (B and E..AD) to rotate between rows within each group of matches
# that were re-ordered in this update. Column A (row index), C (league) and D (date)
# are identical for every row inside a group so they are left untouched.
$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

# Snapshot current (pre-edit) values for every row involved, using Value2 to read the
# raw underlying value (avoids Value returning a wrapped/boxed object).
$snapshot = @{}
$row75 = @{}
foreach ($c in $cols) { $row75[$c] = $ws.Range("${c}75").Value2 }
$snapshot[75] = $row75
$row76 = @{}
foreach ($c in $cols) { $row76[$c] = $ws.Range("${c}76").Value2 }
$snapshot[76] = $row76
$row77 = @{}
foreach ($c in $cols) { $row77[$c] = $ws.Range("${c}77").Value2 }
$snapshot[77] = $row77
$row85 = @{}
foreach ($c in $cols) { $row85[$c] = $ws.Range("${c}85").Value2 }
$snapshot[85] = $row85
$row88 = @{}
foreach ($c in $cols) { $row88[$c] = $ws.Range("${c}88").Value2 }
$snapshot[88] = $row88
$row119 = @{}
foreach ($c in $cols) { $row119[$c] = $ws.Range("${c}119").Value2 }
$snapshot[119] = $row119
$row120 = @{}
foreach ($c in $cols) { $row120[$c] = $ws.Range("${c}120").Value2 }
$snapshot[120] = $row120
$row121 = @{}
foreach ($c in $cols) { $row121[$c] = $ws.Range("${c}121").Value2 }
$snapshot[121] = $row121
$row302 = @{}
foreach ($c in $cols) { $row302[$c] = $ws.Range("${c}302").Value2 }
$snapshot[302] = $row302
$row303 = @{}
foreach ($c in $cols) { $row303[$c] = $ws.Range("${c}303").Value2 }
$snapshot[303] = $row303
$row304 = @{}
foreach ($c in $cols) { $row304[$c] = $ws.Range("${c}304").Value2 }
$snapshot[304] = $row304
$row305 = @{}
foreach ($c in $cols) { $row305[$c] = $ws.Range("${c}305").Value2 }
$snapshot[305] = $row305
$row306 = @{}
foreach ($c in $cols) { $row306[$c] = $ws.Range("${c}306").Value2 }
$snapshot[306] = $row306

# Apply the rotated values to destination rows (mapping: destination row -> source row
# whose pre-edit data it now receives)
foreach ($c in $cols) { $ws.Range("${c}75").Value2 = $snapshot[77][$c] }
foreach ($c in $cols) { $ws.Range("${c}76").Value2 = $snapshot[75][$c] }
foreach ($c in $cols) { $ws.Range("${c}77").Value2 = $snapshot[76][$c] }
foreach ($c in $cols) { $ws.Range("${c}85").Value2 = $snapshot[88][$c] }
foreach ($c in $cols) { $ws.Range("${c}88").Value2 = $snapshot[85][$c] }
foreach ($c in $cols) { $ws.Range("${c}119").Value2 = $snapshot[121][$c] }
foreach ($c in $cols) { $ws.Range("${c}120").Value2 = $snapshot[119][$c] }
foreach ($c in $cols) { $ws.Range("${c}121").Value2 = $snapshot[120][$c] }
foreach ($c in $cols) { $ws.Range("${c}302").Value2 = $snapshot[303][$c] }
foreach ($c in $cols) { $ws.Range("${c}303").Value2 = $snapshot[304][$c] }
foreach ($c in $cols) { $ws.Range("${c}304").Value2 = $snapshot[305][$c] }
foreach ($c in $cols) { $ws.Range("${c}305").Value2 = $snapshot[306][$c] }
foreach ($c in $cols) { $ws.Range("${c}306").Value2 = $snapshot[302][$c] }
